$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round Q4/R4 to whole numbers (integer values)
$ws.Range("Q4").Value = 485109
$ws.Range("R4").Value = 6406776

# Clear the Starttid (Z4) and Sluttid (AB4) cells entirely
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
